# "Add files via upload" — 신승민 adds a new task row describing the
# plan-recommendation algorithm work, and the active sheet/selection in
# the workbook moves from 황석영's sheet to 신승민's sheet.

$wb = $excel.ActiveWorkbook

# ---- 신승민 sheet (first sheet / tab order 1) ----------------------------
$ws1 = $wb.Worksheets.Item(1)

# New row of task data (row 5 was blank, now filled in).
$ws1.Cells.Item(5, 1).Value = "요금제 추천 알고리즘 구현"
$ws1.Cells.Item(5, 2).Value = "요금제 추천 구현"
$ws1.Cells.Item(5, 3).Value = 43611
$ws1.Cells.Item(5, 4).Value = 43612
$ws1.Cells.Item(5, 5).Value = "입력받은 값에 따라 요금제를 추천하는 알고리즘 구현 ( Plan.java , recomPlan.java 구현)"

# ---- 황석영 sheet (second sheet / tab order 2) ---------------------------
# Leave its data alone, just move the selection/cursor.
$ws2 = $wb.Worksheets.Item(2)
[void]$ws2.Activate()
[void]$ws2.Range("B4").Select()

# ---- Make 신승민's sheet the active tab/selection again -----------------
[void]$ws1.Activate()
[void]$ws1.Range("F5").Select()
